$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3: change Condition to "contains" and Search_Term to "a"
$ws.Range("E3").Value = "a"
$ws.Range("D3").Value = "contains"

# Add a new row 4 with a City / ends with / s rule
$ws.Range("C4").Value = "City"
$ws.Range("D4").Value = "ends with"
$ws.Range("E4").Value = "s"

# Leave selection on E5, matching the saved cursor position
$ws.Range("E5").Select()

$wb.Save()
